# Orders workbook update: record lactose/nut-free flags for the newest
# order (Candy Canes - Green, row 10) now that it has been appended to
# inventory with a quantity, and roll the flags that used to live on the
# now-removed row 11 order forward/back across the other candy rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Creme Eggs - Large: now flagged has_lactose = Y (was N)
$ws.Range("R7").Value = "Y"

# Row 10 - Candy Canes - Green: newly given has_lactose / has_nuts values
$ws.Range("R10").Value = "N"
$ws.Range("S10").Value = "N"

# Row 11 - Easter Bunny with Beach Attire: has_lactose/has_nuts cleared,
# these fields no longer apply to this (non-candy) order
$ws.Range("R11").Value = ""
$ws.Range("S11").Value = ""

# Row 14 - Candy Canes - Red: has_lactose flipped to N (was Y)
$ws.Range("R14").Value = "N"

# Row 15 - Sea Salt Pumpkin Caramel Toffee: has_lactose flipped to Y (was N)
$ws.Range("R15").Value = "Y"

# Reflect the updated working selection/scroll position from the edit session
$ws.Range("T14").Select()
